# Kientrucdetai_26_8.pptx -- "Chinh sua ngay 31/8"
#
# 1) Every "Date Placeholder" shape (datetimeFigureOut field) on the slide
#    master and on all 11 slide layouts shows the cached date 8/29/2010;
#    bump it to 8/30/2010.
# 2) On slide 1, the "Flowchart: Magnetic Disk 26" shape moves down by
#    12pt (152400 EMU / 12700 EMU-per-pt): Top 234pt -> 246pt
#    (y 2971800 -> 3124200 EMU).

$p = $ppt.ActivePresentation

$oldDate = "8/29/2010"
$newDate = "8/30/2010"

function Update-DateField($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# --- Slide master date placeholder ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateField $master.Shapes.Item($i)
}

# --- Every custom (slide) layout's date placeholder ---
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateField $layout.Shapes.Item($i)
    }
}

# --- Reposition "Flowchart: Magnetic Disk 26" on slide 1 ---
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Name -eq "Flowchart: Magnetic Disk 26") {
        $shp.Top = 246
    }
}
